$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E8").Select()
$ws.Range("E8").Value = "GIT UPDATE"
